$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 752, shifting existing rows 752-793 down to 753-794.
$ws.Rows.Item(752).Insert()

# Column A ("2026/02/04") looks like a date, so force text formatting before
# assigning the value, then restore the default style so the cell matches the
# unstyled data rows around it.
$ws.Range("A752").NumberFormat = "@"
$ws.Range("A752").Value = "2026/02/04"
$ws.Range("A752").Style = "Normal"

$ws.Range("B752").Value = "水"
$ws.Range("C752").Value = 7
$ws.Range("D752").Value = 27
